# Add a new "Swiss" worksheet (Switzerland market test data) to the
# workbook, cloned from the existing "Czech" sheet template, and leave it
# as the active sheet/tab.

$wb = $excel.ActiveWorkbook
$czech = $wb.Worksheets.Item("Czech")

# Czech sheet is no longer the active tab once Swiss is added; its saved
# selection becomes a "select all" (entire grid) selection.
$czech.Cells.Select() | Out-Null

# Clone the Czech sheet (same columns/styles/layout) and place it right
# after Czech, at the end of the tab strip.
$czech.Copy([System.Reflection.Missing]::Value, $czech)
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Fill in the Switzerland-specific test data.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2653/T2654/T2655"

# Make the new sheet the active tab with its own saved selection.
$swiss.Activate() | Out-Null
$swiss.Range("B15").Select() | Out-Null
